# Automated update by travis-ci: `make clean download import process`
#
# Adds a new candidate-totals row for "Michael Houston"
# (oakland-2024-03-05 election) to the "Totals Per Candidates" sheet,
# inserted directly above the existing "Bruce Quan" row (which, along
# with every row below it, shifts down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at sheet row 137 - this shifts the old row 137
# (Bruce Quan) and everything below it down by one row.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row with the new candidate's record.
$ws.Cells.Item(137, 1).Value = "michael-houston.json"
$ws.Cells.Item(137, 2).Value = "oakland-2024-03-05"
$ws.Cells.Item(137, 3).Value = "Michael Houston"
$ws.Cells.Item(137, 4).Value = 17045
$ws.Cells.Item(137, 5).Value = 17045
$ws.Cells.Item(137, 6).Value = 17045
$ws.Cells.Item(137, 7).Formula = "=E137-D137"
$ws.Cells.Item(137, 8).Formula = "=F137-D137"

# Keep the formula cells unstyled (matching the rest of the sheet) - the
# row Insert() otherwise carries a stray number format onto them.
$ws.Cells.Item(137, 7).Style = "Normal"
$ws.Cells.Item(137, 8).Style = "Normal"
